# Mars classifier list — reclassification edit.
#
# Commit: "Inclusion of Marsupial classifier list as a text file (tab-delimited)"
# The underlying data was round-tripped through a tab-delimited text file, which:
#   - collapsed several compound Diet categories into their simpler parent
#     category (e.g. "Carnivore-insectivore" -> "Insectivore",
#     "Carnivore-omnivore" -> "Carnivore", "Herbivore-omnivore" -> "Herbivore")
#   - replaced punctuation ("/", "-") that doesn't survive a plain
#     tab-delimited round trip with underscores in the remaining composite
#     category labels (e.g. "Aquatic/Terrestrial" -> "Aquatic_Terrestrial",
#     "Herbivore-folivore" -> "Herbivore_folivore")
#   - left the Specimen names (column A) and all other Diet/Locomotion/Clade
#     values untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Diet (column C) / Locomotion (column B) relabeling -------------------

$ws.Range("C3").Value  = "Herbivore"
$ws.Range("C4").Value  = "Insectivore"

$ws.Range("B6").Value  = "Aquatic_Terrestrial"
$ws.Range("C6").Value  = "Carnivore"

$ws.Range("C7").Value  = "Insectivore"

$ws.Range("C12").Value = "Insectivore_nectarivore"

$ws.Range("C13").Value = "Insectivore"
$ws.Range("C15").Value = "Insectivore"

$ws.Range("C16").Value = "Herbivore_folivore"

$ws.Range("C19").Value = "Insectivore"

$ws.Range("C21").Value = "Herbivore_folivore"

$ws.Range("B23").Value = "Arboreal_Scansorial"
$ws.Range("C23").Value = "Insectivore"

$ws.Range("B24").Value = "Arboreal_Scansorial"
$ws.Range("C24").Value = "Insectivore"

$ws.Range("B30").Value = "Terrestrial_Scansorial"
$ws.Range("C30").Value = "Insectivore"

$ws.Range("C33").Value = "Herbivore_folivore"

$ws.Range("C34").Value = "Nectarivore_omnivore"

$ws.Range("C36").Value = "Herbivore_folivore"

$ws.Range("C41").Value = "Insectivore"
$ws.Range("C43").Value = "Insectivore"
$ws.Range("C46").Value = "Insectivore"
$ws.Range("C47").Value = "Insectivore"

$ws.Range("C49").Value = "Herbivore_folivore"
$ws.Range("C50").Value = "Herbivore_folivore"

# --- View state: scroll position + active selection moved to K48 ----------

$excel.ActiveWindow.ScrollRow = 29
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("K48").Select()
